$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: membership counts -- leading apostrophe keeps these stored as literal
# text (shared strings), matching the original cell type, instead of Excel
# auto-converting the numeric-looking input into a Number cell.
$ws.Range("B2").Value = "'800"
$ws.Range("B3").Value = "'320"
$ws.Range("B4").Value = "'100"
$ws.Range("B5").Value = "'136"
$ws.Range("B6").Value = "'500"
$ws.Range("B2:B6").ClearFormats()

# Column C: Encompasses community sites
$ws.Range("C2").Value = 'No, FLASCO does not include community sites, as it is focused on clinical oncology professionals in Florida.'
$ws.Range("C3").Value = 'No, GASCO does not encompass community sites because it is focused on connecting and supporting clinical oncology professionals in Georgia.'
$ws.Range("C4").Value = 'No, IOS primarily focuses on academic institutions for oncology research and care, not community sites.'
$ws.Range("C5").Value = 'No, The IOWA Oncology Society focuses primarily on academic and research institutions, not community sites.'
$ws.Range("C6").Value = 'No, MOASC does not encompass community sites, as it primarily focuses on medical oncology practices rather than community-based facilities.'

# Column D: Influential on state or local policy
$ws.Range("D2").Value = 'No, FLASCO primarily focuses on education and advocacy within the oncology community, rather than influencing state or local policy.'
$ws.Range("D3").Value = 'No, GASCO is primarily focused on education and advocacy for clinical oncologists, rather than on influencing policy.'
$ws.Range("D4").Value = 'No, IOS is primarily a professional organization for oncologists and researchers and does not have a direct influence on state or local policy.'
$ws.Range("D5").Value = 'No, 
The IOWA Oncology Society is not influential on state or local policy because its focus is primarily on promoting education, research, and professional development within the field of oncology rather than advocating for policy changes.'
$ws.Range("D6").Value = 'No, MOASC is primarily focused on education and networking for medical oncologists in Southern California, not on influencing policy.'

# Column E: Engagement opportunity with leadership
$ws.Range("E2").Value = 'Yes, FLASCO provides engagement opportunity with leadership, as they offer various networking events and conferences where members can interact with the organization''s leadership team.'
$ws.Range("E3").Value = 'Yes, GASCO provides engagement opportunity with leadership, they offer various networking events and leadership development programs.'
$ws.Range("E4").Value = 'Yes, IOS provides engagement opportunities with leadership. IOS offers various networking events, conferences, and leadership development programs where members can interact and engage with key leaders in the oncology field.'
$ws.Range("E5").Value = 'Yes, 
The IOWA Oncology Society does provide engagement opportunities with leadership through networking events, conferences, and mentorship programs.'
$ws.Range("E6").Value = 'Yes, the organization provides engagement opportunities with leadership through various events, meetings, and networking opportunities where members can interact with top leaders in the field of medical oncology.'

# Column F: Support for clinical trial recruitment
$ws.Range("F2").Value = 'Yes, FLASCO does provide support for clinical trial recruitment through education and resources for members to refer patients to appropriate trials.'
$ws.Range("F3").Value = 'No, GASCO does not provide support for clinical trial recruitment. GASCO''s primary focus is on providing education, resources, and advocacy for oncology professionals in Georgia.'
$ws.Range("F4").Value = 'No, IOS does not provide support for clinical trial recruitment. Clinical trial recruitment typically falls under the purview of academic institutions and research organizations, rather than professional societies like IOS.'
$ws.Range("F5").Value = 'Yes, 
The IOWA Oncology Society works closely with research institutions and pharmaceutical companies to promote clinical trial recruitment among its members.'
$ws.Range("F6").Value = 'No, MOASC does not typically provide support for clinical trial recruitment. They primarily focus on education and networking for medical oncology professionals in Southern California.'

# Column G: Engagement opportunity with payors
$ws.Range("G2").Value = 'No, FLASCO does not provide engagement opportunity with payors. FLASCO focuses on providing education, resources, and support for clinical oncologists, rather than direct interactions with payors.'
$ws.Range("G3").Value = 'No, GASCO does not have direct engagement opportunities with payors. GASCO primarily focuses on providing education, advocacy, and networking opportunities for oncology professionals in Georgia.'
$ws.Range("G4").Value = 'No, IOS does not provide engagement opportunity with payors. IOS is a professional organization focused on oncology practice in Indiana, rather than a platform for engagement with payors.'
$ws.Range("G5").Value = 'No, The focus of IOWA Oncology Society is primarily on oncology education and research, not payer engagement.'
$ws.Range("G6").Value = 'No, MOASC does not provide engagement opportunity with payors. The association focuses on providing education and networking opportunities for medical oncologists in the Southern California region.'

# Column H: Area experts on board
$ws.Range("H2").Value = 'No, FLASCO does not include area experts on its board. FLASCO is a professional organization for clinical oncology professionals and not specifically for area experts.'
$ws.Range("H3").Value = 'Yes, GASCO includes area experts on its board. GASCO is an organization dedicated to clinical oncology, so it would be vital to have experts on the board to provide insight and guidance.'
$ws.Range("H4").Value = 'No, , IOS does not have area experts on its board.'
$ws.Range("H5").Value = 'Yes, the IOWA Oncology Society includes area experts on its board to ensure high-quality leadership and decision-making.'
$ws.Range("H6").Value = 'No, 
The Medical Oncology Association of Southern California (MOASC) does not have area experts on its board, as the board members primarily consist of oncologists and healthcare professionals in the field.'

# Column I: Therapeutic research collaborations
$ws.Range("I2").Value = 'Yes, FLASCO is involved in therapeutic research collaborations. FLASCO actively engages in collaborations with various organizations and institutions to advance research in oncology.'
$ws.Range("I3").Value = 'Yes, GASCO is involved in therapeutic research collaborations. GASCO actively partners with academic institutions and pharmaceutical companies to conduct clinical trials and advance cancer research.'
$ws.Range("I4").Value = 'Yes, IOS is involved in therapeutic research collaborations. IOS actively participates in collaborative research efforts to advance oncology care and treatment options.'
$ws.Range("I5").Value = 'No, IOWA Oncology Society primarily focuses on education and networking opportunities for oncology professionals.'
$ws.Range("I6").Value = 'Yes, MOASC is involved in therapeutic research collaborations. The association actively partners with pharmaceutical companies, academic institutions, and other organizations to advance research in the field of oncology.'

# Column J: Top therapeutic area experts on board
$ws.Range("J2").Value = 'No, FLASCO does not include top therapeutic area experts on its board. The organization primarily consists of clinical oncologists and healthcare professionals who specialize in oncology management and research.'
$ws.Range("J3").Value = 'Yes, GASCO includes top therapeutic area experts on its board. GASCO is a professional organization for clinical oncologists, so it is likely that its board members are experts in various areas of oncology.'
$ws.Range("J4").Value = 'No, While IOS likely includes experts in oncology on its board, it cannot be confirmed whether they are all top therapeutic area experts.'
$ws.Range("J5").Value = 'Yes, the IOWA Oncology Society includes top therapeutic area experts on its board, as they are dedicated to advancing oncology care and treatment.'
$ws.Range("J6").Value = 'Yes, MOASC includes top therapeutic area experts on its board. The organization is comprised of medical oncologists and researchers who specialize in various areas of cancer treatment, providing a wealth of expertise and knowledge to guide its decision-making and initiatives.'

# A couple of the new justifications embed a literal newline inside the cell text,
# which makes Excel auto-expand those rows height when the value is assigned.
# Re-run AutoFit on just those rows so the stored row height matches the untouched original.
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()